$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "custom accuracy": reduce the stored precision of the row-5 measurement
# cells (B5:AH5) to 2 decimal places.
$roundedRow5 = [ordered]@{
    "B5"  = 19.22
    "C5"  = 14.09
    "D5"  = 1.2
    "E5"  = 41.75
    "F5"  = 34.13
    "G5"  = 15.12
    "H5"  = 58.92
    "I5"  = 23.27
    "J5"  = 10.3
    "K5"  = 15.22
    "L5"  = 16.76
    "M5"  = 17.64
    "N5"  = 4.83
    "O5"  = 15.04
    "P5"  = 21.37
    "Q5"  = 12.72
    "R5"  = 0.84
    "S5"  = 0.8
    "T5"  = 221.59
    "U5"  = 42.07
    "V5"  = 13.88
    "W5"  = 28.22
    "X5"  = 14.77
    "Y5"  = 2.29
    "Z5"  = 28.6
    "AA5" = 12.26
    "AB5" = 10.91
    "AC5" = 12.82
    "AD5" = 17.54
    "AE5" = 0.5600000000000001
    "AF5" = 53.51
    "AG5" = 7.79
    "AH5" = 17.35
}

foreach ($addr in $roundedRow5.Keys) {
    $ws.Range($addr).Value = $roundedRow5[$addr]
}

# "데이터 1000개" (trim the dataset to 1000 rows): drop the extra
# simulation timestep that used to live in row 6.
$ws.Rows.Item(6).Delete()
